$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.030.55"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.562.53"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'208.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("D11").Value = "'0.0854"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "1.785.30"
$ws.Range("D13").Value = "1.561.40"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "'0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "27.033.75"
$ws.Range("D17").Value = "'61.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "0.0₃0708"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "'216.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").Value = "'7.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").Value = "'9.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").Value = "'1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").Value = "'153.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").Value = "'15.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("E31").Value = "  +3.16%  "
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("D34").Value = "1.431.90"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("E36").Value = "  +8.17%  "
$ws.Range("D37").Value = "'2.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("D39").Value = "'0.534"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.58%  "
$ws.Range("D40").Value = "'5.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.53%  "
$ws.Range("D41").Value = "'0.809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "'64.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("D47").Value = "1.700.18"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "'87.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("E49").Value = "  +4.15%  "
$ws.Range("D51").Value = "'0.0958"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.35%  "
